{"js": "// Blazing Goddess review: add a \"Meta description\" paragraph right after\n// the title, drop the old duplicate \"Play Blazing Goddess...\" heading that\n// was sitting near the end of the document, and turn the trailing italic\n// paragraph into the AI image-generation prompt.\n\nconst boldLabel = \"Meta description\";\nconst restOfSentence =\n  \": Find out all about the gameplay mechanics, graphics, symbols, and ways to win in Blazing Goddess. Play for free and read the game review now.\";\nconst promptText =\n  \"Prompt: Create a feature image for Blazing Goddess that reflects the game's exciting and adventurous vibe. The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a flaming torch in one hand and a coconut cocktail in the other. The backdrop should be a tropical paradise with palm trees, sand, and sea. The Blazing Goddess symbol should be prominently displayed in the background, with flames erupting from it. The overall aesthetic should be bright, colorful, and energetic to capture the excitement of playing the Blazing Goddess slot machine.\";\n\n// --- 1) Insert a new \"Meta description\" paragraph right after the H1 title ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst metaParagraph = titleParagraph.insertParagraph(\n  boldLabel + restOfSentence,\n  \"After\"\n);\nmetaParagraph.style = \"Normal\";\nawait context.sync();\n\n// Bold only the \"Meta description\" label, leave the rest of the sentence regular.\nconst boldMatches = metaParagraph.search(boldLabel, { matchCase: true });\nboldMatches.load(\"items\");\nawait context.sync();\nboldMatches.items[0].font.bold = true;\nawait context.sync();\n\n// --- 2) Remove the bold \"Play Blazing Goddess for Free | Review of Slot Game\"\n//        paragraph that duplicated the title near the bottom of the document ---\nconst allParagraphs = context.document.body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nconst count = allParagraphs.items.length;\nconst duplicateTitleParagraph = allParagraphs.items[count - 2];\nduplicateTitleParagraph.delete();\nawait context.sync();\n\n// --- 3) Replace the final (italic) paragraph's text with the new image prompt ---\nconst remainingParagraphs = context.document.body.paragraphs;\nremainingParagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph =\n  remainingParagraphs.items[remainingParagraphs.items.length - 1];\nlastParagraph.insertText(promptText, \"Replace\");\nawait context.sync();\n", "ps1": "# Blazing Goddess review: add a \"Meta description\" paragraph right after\n# the title, drop the old duplicate \"Play Blazing Goddess...\" heading that\n# was sitting near the end of the document, and turn the trailing italic\n# paragraph into the AI image-generation prompt.\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert a new \"Meta description\" paragraph right after the H1 title ---\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Style = \"Normal\"\n\n$boldText = \"Meta description\"\n$restText = \": Find out all about the gameplay mechanics, graphics, symbols, and ways to win in Blazing Goddess. Play for free and read the game review now.\"\n\n$metaRange = $metaPara.Range\n$metaRange.Collapse(1)\n$metaRange.InsertAfter($boldText + $restText)\n\n# Bold only the \"Meta description\" label, leave the rest of the sentence regular.\n$metaStart = $metaPara.Range.Start\n$labelRange = $d.Range($metaStart, $metaStart + $boldText.Length)\n$labelRange.Bold = 1\n\n# --- 2) Remove the bold \"Play Blazing Goddess for Free | Review of Slot Game\"\n#        paragraph that duplicated the title near the bottom of the document ---\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n# --- 3) Replace the final (italic) paragraph's text with the new image prompt ---\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastRangeRaw = $lastPara.Range\n$lastRange = $d.Range($lastRangeRaw.Start, $lastRangeRaw.End)\n$lastRange.Text = \"Prompt: Create a feature image for Blazing Goddess that reflects the game's exciting and adventurous vibe. The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a flaming torch in one hand and a coconut cocktail in the other. The backdrop should be a tropical paradise with palm trees, sand, and sea. The Blazing Goddess symbol should be prominently displayed in the background, with flames erupting from it. The overall aesthetic should be bright, colorful, and energetic to capture the excitement of playing the Blazing Goddess slot machine.\"\n"}
